$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: propagate cell formatting (styles) to new row positions ---
# Capture the two less-common patterns (rows 6 and 7) before they get overwritten.
$ws.Range("A6:E6").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)
$ws.Range("A7:E7").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)
$ws.Range("A7:E7").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)
$ws.Range("A5:E5").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)

# --- Step 2: set cell values for the job-step rows ---
# Row 6
$ws.Range("A6").Value = "CoachingInactivations"
$ws.Range("B6").Value = "Inactivations"
$ws.Range("C6").Value = "Inactivations logs"
$ws.Range("D6").Value = "\\vrivscors01\BCC Scorecards\Coaching\Inactivations\eCL_Coaching_<MMDDCCYY>.csv`n\\vrivscors01\BCC Scorecards\Coaching\Inactivations\eCL_Warning_<MMDDCCYY>.csv"
$ws.Range("E6").Value = "EC.Coaching_Log.StatusID = 2 or`nEC.Warning_Log.StatusID = 2`nemail notification sent to john;`nlog file generated to \\vrivscors01\BCC Scorecards\Coaching\Inactivations\Processed"

# Row 8
$ws.Range("A8").Value = "CoachingNotifications"
$ws.Range("B8").Value = "Notifications"
$ws.Range("C8").Value = "Sends emails to coaching log recipients "
$ws.Range("D8").Value = "EC.Coaching_Log.EmailSent = 0 or null"
$ws.Range("E8").Value = "EC.Coaching_Log.EmailSent = 1`nemail notification sent to recipient"

# Row 10
$ws.Range("A10").Value = "CoachingOutliersLoad"
$ws.Range("B10").Value = "OutliersLoad"
$ws.Range("C10").Value = "Imports OMR information to create coaching logs"
$ws.Range("D10").Value = "\\vrivscors01\BCC Scorecards\Coaching\Outliers\Test\eCl_Outlier_Feed_<ReportCode><CCYYMMDD>.csv"
$ws.Range("E10").Value = "EC.Outlier_Coaching_Stage`n    EC.Outlier_Coacing_Rejected`n    EC.Outlier_Coaching_Fact`n        EC.Coaching_Log`n        EC.Coaching_Log_Reason"

# Row 12
$ws.Range("A12").Value = "CoachingQualityLoad"
$ws.Range("B12").Value = "IQSLoad"
$ws.Range("C12").Value = "Imports IQS quality call information to create coaching logs"
$ws.Range("D12").Value = "\\vrivscors01\BCC Scorecards\Coaching\Apps\Encryption\Encrypt_out\eCL_IQS_Scorecard_<CCYYMMDD>.csv.zip.encrypt"
$ws.Range("E12").Value = "EC.Quality_Coaching_Stage`n    EC.Quality_Coacing_Rejected`n    EC.Quality_Coaching_Fact`n        EC.Coaching_Log`n        EC.Coaching_Log_Reason"

# Row 14
$ws.Range("A14").Value = "CoachingReminders"
$ws.Range("B14").Value = "Reminders"
$ws.Range("C14").Value = "Sends email reminders when certain logs have not been coached"
$ws.Range("D14").Value = "selection determined by log source/reason/sub-reason/value"
$ws.Range("E14").Value = "EC.Coaching_Log"

# Row 16
$ws.Range("A16").Value = "CoachingSurveyGenerate"
$ws.Range("B16").Value = "SurveyGenerate"
$ws.Range("C16").Value = "Provides for selecting those logs which will be included in survey"
$ws.Range("D16").Value = "EC.Coaching_Log"
$ws.Range("E16").Value = "EC.Survey_Response_Detail"

# Row 18
$ws.Range("A18").Value = "CoachingSurveyNotifications"
$ws.Range("B18").Value = "SurveyNotifications"
$ws.Range("C18").Value = "Sends email notification regarding survey"
$ws.Range("D18").Value = "EC.Coaching_Log.SurveySent = 0"
$ws.Range("E18").Value = "EC.Coaching_Log.SurveySent = 1`nemail notification sent to recipient"

# Row 20
$ws.Range("A20").Value = "CoachingWHLoad"
$ws.Range("B20").Value = "WHFileLoad"
$ws.Range("C20").Value = "Imports information to create warning logs"
$ws.Range("D20").Value = "\\vrivscors01\BCC Scorecards\Coaching\WH\<SiteLocation>Warnings.csv"
$ws.Range("E20").Value = "EC.Warning_History_Stage`n    EC.Warning_History_Rejected`n    EC.Warning_History_Fact`n        EC.Warning_Log`n        EC.Warning_Log_Reason"

$excel.CutCopyMode = 0
